$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: serial_number -> serialNumber
$ws.Range("A1").Value = "serialNumber"

# Update serial numbers in A2:A6
$ws.Range("A2").Value = "SN18123"
$ws.Range("A3").Value = "SN18124"
$ws.Range("A4").Value = "SN18125"
$ws.Range("A5").Value = "SN18126"
$ws.Range("A6").Value = "SN18127"

# Update selection to F15
$ws.Range("F15").Select()
